$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final sorted dataset (alphabetical by name), rows 2..12
$data = @(
    @("Ja",   "Alexander Götz von olenhusen", 49,  0, "hoodie_zip man L (49€)"),
    @("Ja",   "Annika Nolte", 49, 0, "hoodie_zip woman M (49€)"),
    @("Nein", "Chris Bierfreund", 167, 1, "jersey dark man M (34€), jersey light man M (34€), jersey_long black man M (38€), short multisport S (22€), hoodie_zip man M (49€)"),
    @("Nein", "Elena Gugganig", 80, 1, "jersey dark woman L (34€), jersey light woman L (34€), short woman M (22€)"),
    @("Ja",   "Hannes Drobek", 109, 0, "jersey_long black man L (38€), tank dark man L (29€), hoodie_no_zip man L (42€)"),
    @("Nein", "Ingo Schmiegel", 209, 1, "jersey dark man L (34€), jersey light man L (34€), jersey_long black man L (38€), short multisport L (22€), hoodie_no_zip man L (42€), hoodie_zip man L (49€)"),
    @("Ja",   "Julia Busse", 100, 0, "tank light woman L (29€), short multisport L (22€), hoodie_zip woman M (49€)"),
    @("Ja",   "Katharina Wilbrandt", 22, 0, "short multisport XL (22€)"),
    @("Ja",   "Matthis Heimberg", 109, 0, "jersey_long black man L (38€), tank dark man L (29€), hoodie_no_zip man L (42€)"),
    @("Nein", "Nils Hoffmann", 80, 1, "jersey dark man S (34€), jersey light man S (34€), short multisport M (22€)"),
    @("Nein", "Patrick Jung", 80, 1, "jersey dark man L (34€), jersey light man L (34€), short long L (22€)")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $row++
}
